# Update countries & provincias Spain
# Applies the latest COVID snapshot update to the "Pais" worksheet:
#  - refreshes counters for a handful of countries
#  - re-inserts "Alemania", "Estado de Palestina" and "Eslovenia" earlier in
#    the country list (their rows move up while carrying their refreshed
#    figures; the countries they displace simply shift down one row,
#    keeping their own, unchanged figures)
#  - bumps the "last updated" timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 16 de Octubre de 2020 a las 15:49"

# Each entry: row, country name (col A), then B,C,D,E,F,G,H
$rows = @(
    @(4,  "Estados Unidos", 8222371, 6056, 5320751, 2678789, 0, 114, 222831),
    @(5,  "India",          7376783, 11274, 6458305, 806199, 0, 133, 112279),
    @(18, "Irak",            420303, 3501,  353962,  56199,  0, 56,  10142),
    @(22, "Alemania",        352401, 3585,  284600,  57976,  0, 15,  9825),
    @(23, "Filipinas",       351750, 3139,  294865,  50354,  0, 34,  6531),
    @(25, "Arabia Saudita",  341495, 433,   327795,  8556,   0, 17,  5144),
    @(49, "Portugal",        95902,  2608,  56066,   37687,  0, 21,  2149),
    @(60, "Uzbekistan",      62588,  310,   59624,   2444,   0, 2,   520),
    @(71, "Estado de Palestina", 46434, 334, 39921, 6111,    0, 1,   402),
    @(72, "Irlanda",          46429, 0,   23364,   21227,  0, 0,   1838),
    @(77, "Serbia",           35719, 265, 31536,   3411,   0, 2,   772),
    @(96, "Noruega",          16201, 65,  11863,   4060,   0, 0,   278),
    @(104,"Eslovenia",        11517, 834, 5924,    5413,   0, 4,   180),
    @(105,"Guinea",           11327, 72,  10380,   877,    0, 0,   70),
    @(106,"Maldivas",         11113, 0,   9931,    1147,   0, 0,   35),
    @(107,"Consejo Danes para los Refugiados", 10935, 0, 10306, 348, 0, 0, 281),
    @(109,"Tayikistan",       10414, 40,  9393,    941,    0, 0,   80)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}
